$d = $word.ActiveDocument

$replacements = @(
    @("63×79=", "89×86="),
    @("87×45=", "56×79="),
    @("46×12=", "78×99="),
    @("89×40=", "65×42="),
    @("38×90=", "64×71="),
    @("83×31=", "17×26="),
    @("53×72=", "48×11="),
    @("36×66=", "24×39="),
    @("76×26=", "29×41="),
    @("70×82=", "59×89="),
    @("34×77=", "83×19="),
    @("21×90=", "75×72="),
    @("55×31=", "52×46="),
    @("26×25=", "25×28="),
    @("94×91=", "85×98="),
    @("78×75=", "60×18="),
    @("39×28=", "70×34="),
    @("58×92=", "89×58="),
    @("37×92=", "56×85="),
    @("41×93=", "74×23="),
    @("96×70=", "30×71="),
    @("95×78=", "60×40="),
    @("83×15=", "34×25="),
    @("13×48=", "11×70="),
    @("51×98=", "88×67=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
